$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the test case IDs (GESACT -> GESTACT)
$ws.Range("A2").Value = "CP_GESTACT_001"
$ws.Range("A3").Value = "CP_GESTACT_002"
$ws.Range("A4").Value = "CP_GESTACT_003"

# Make header row taller
$ws.Rows.Item(1).RowHeight = 30

# Update view: scroll back to A1, select G4
$ws.Range("G4").Select()
